# "Editing users on IC" — adds new IC user-update test rows on the
# "ic_login++" sheet, adds corresponding scenario rows on the "IC" sheet,
# and updates which sheet/cell is active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "ic_login++" sheet — replace the old test user with the new one and
#    extend the table from 3 data rows (TCID 4-6) to 10 data rows
#    (TCID 4-13), each carrying a Username/Password mailto hyperlink.
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("ic_login++")

$urlValue = "https://staging-jdgroup-m23.vaimo.net/T5sjY7drHkyB6Z4n/admin/index/index/key/4cee16108e07c6904ab12f33e9f10b9b1b1fadb7c06faa4c94fbd03a1d1018ff/"
$userValue = "LauraPLittle@teleworm.us"
$passValue = "Password@123"

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $tcid = 4 + $i

    $wsLogin.Cells.Item($row, 1).Value = $tcid
    $wsLogin.Cells.Item($row, 2).Value = 1
    $wsLogin.Cells.Item($row, 3).Value = $urlValue
    $wsLogin.Cells.Item($row, 4).Value = $userValue
    $wsLogin.Cells.Item($row, 5).Value = $passValue
    $wsLogin.Cells.Item($row, 4).Style = "Hyperlink"
    $wsLogin.Cells.Item($row, 5).Style = "Hyperlink"
}
# Row 2's E column keeps the default (non-hyperlink) style, matching the
# authored workbook.
$wsLogin.Cells.Item(2, 5).Style = "Normal"

# Hyperlinks: column E always links to the mailto password; column D links
# to the mailto username from row 2 onward (row 2/row 3-4 combined range,
# then one per row for the newly-added rows).
$wsLogin.Hyperlinks.Add($wsLogin.Range("D2"), "mailto:$userValue") | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("D3:D4"), "mailto:$userValue", "", "", $userValue) | Out-Null
for ($i = 5; $i -le 11; $i++) {
    $wsLogin.Hyperlinks.Add($wsLogin.Range("E$i"), "mailto:$passValue") | Out-Null
    $wsLogin.Hyperlinks.Add($wsLogin.Range("D$i"), "mailto:$userValue") | Out-Null
}

$wsLogin.Activate()
$wsLogin.Range("A3:A11").Select()

# ---------------------------------------------------------------------
# 2) "IC" sheet — extend the "user update in IC" scenario rows (previously
#    3 generic rows) into the full set of 8 specific user-update scenarios.
# ---------------------------------------------------------------------
$wsIC = $wb.Worksheets.Item("IC")

$scenarioAction = "user update in IC"
$scenarios = @(
    "user update in IC  All updates",
    "user update in IC first name",
    "user update in IC last name",
    "user update in IC vat/tax",
    "user update in IC  email",
    "user update in IC password",
    "user update in IC  billing address",
    "user update in IC  shipping address"
)

for ($i = 0; $i -lt $scenarios.Length; $i++) {
    $row = 5 + $i
    $tcid = 4 + $i

    $wsIC.Cells.Item($row, 1).Value = $tcid
    $wsIC.Cells.Item($row, 2).Value = $scenarioAction
    $wsIC.Cells.Item($row, 3).Value = $scenarios[$i]
    $wsIC.Cells.Item($row, 4).Value = "no"
    $wsIC.Cells.Item($row, 5).Value = "Leverch"
    $wsIC.Cells.Item($row, 6).Value = "ic_login"
    $wsIC.Cells.Item($row, 7).Value = "ICUpdateUser"
}

# Row 7 (TCID 6 / "user update in IC last name") keeps its original
# "yes" flag in column D, unlike the others which were switched to "no".
$wsIC.Cells.Item(7, 4).Value = "yes"

$wsIC.Activate()
$wsIC.Range("D12").Select()

# ---------------------------------------------------------------------
# 3) "ICUpdateUser++" sheet loses the active-tab / top-left-cell state
#    (IC becomes the active sheet instead) but keeps its own selection.
# ---------------------------------------------------------------------
$wsUpdUser = $wb.Worksheets.Item("ICUpdateUser++")
$wsUpdUser.Activate()
$wsUpdUser.Range("C2").Select()

# Re-activate "IC" last so it ends up the active sheet / tab in the
# saved workbook (matches activeTab pointing at "IC").
$wsIC.Activate()
$wsIC.Range("D12").Select()
